# Updated cryptos list on Wed Jan 17 19:50:39 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "42.643.38"
$ws.Range("E2").Value = "  -1.06%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.537.55"
$ws.Range("E3").Value = "  -1.26%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.04%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'308.44"
$ws.Range("E5").Value = "  -2.04%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'100.62"
$ws.Range("E6").Value = "  +4.07%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -1.19%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.16%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -2.01%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "'35.90"
$ws.Range("E10").Value = "  +1.21%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -1.13%  "

# Row 12 - Polkadot
$ws.Range("D12").Value = "'7.36"
$ws.Range("E12").Value = "  -1.25%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +0.06%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "2.934.38"
$ws.Range("E14").Value = "  -0.85%  "

# Row 15 - Chainlink
$ws.Range("D15").Value = "'15.91"
$ws.Range("E15").Value = "  +5.64%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "2.546.33"
$ws.Range("E16").Value = "  -2.18%  "

# Row 17 - Polygon
$ws.Range("E17").Value = "  -3.19%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "42.653.93"
$ws.Range("E18").Value = "  -1.10%  "

# Row 19 - Uniswap
$ws.Range("D19").Value = "'6.76"
$ws.Range("E19").Value = "  -1.04%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  -0.75%  "

# Row 21 - InternetComputer(DFINITY)
$ws.Range("D21").Value = "'12.24"
$ws.Range("E21").Value = "  -3.03%  "

# Row 22 - Litecoin
$ws.Range("D22").Value = "'69.43"
$ws.Range("E22").Value = "  +0.20%  "

# Row 23 - BitcoinCash
$ws.Range("D23").Value = "'243.83"
$ws.Range("E23").Value = "  -3.80%  "

# Row 24 - PancakeSwap
$ws.Range("D24").Value = "'2.90"
$ws.Range("E24").Value = "  -2.15%  "

# Row 25 - ImmutableX
$ws.Range("E25").Value = "  -1.39%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.05%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "'26.06"
$ws.Range("E27").Value = "  -3.14%  "

# Row 28 - Toncoin
$ws.Range("E28").Value = "  -3.82%  "

# Row 29 - InjectiveProtocol
$ws.Range("D29").Value = "'39.31"
$ws.Range("E29").Value = "  -2.24%  "

# Row 30 - Cosmos
$ws.Range("E30").Value = "  -1.01%  "

# Row 31 - Filecoin
$ws.Range("D31").Value = "'5.78"
$ws.Range("E31").Value = "  -0.84%  "

# Row 32 - Monero
$ws.Range("D32").Value = "'156.09"
$ws.Range("E32").Value = "  +0.83%  "

# Row 33 - ApeXProtocol
$ws.Range("D33").Value = "'2.72"
$ws.Range("E33").Value = "  +10.91%  "

# Row 34 - Hedera
$ws.Range("D34").Value = "'0.0793"
$ws.Range("E34").Value = "  -1.60%  "

# Row 36 - Row 36 <- Celestia (swap)
$ws.Range("B36").Value = "Celestia"
$ws.Range("C36").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D36").Value = "'18.31"
$ws.Range("E36").Value = "  -3.58%  "

# Row 37 - Row 37 <- ARBITRUM (swap)
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "'2.02"
$ws.Range("E37").Value = "  -4.50%  "

# Row 38 - LidoDAOToken
$ws.Range("E38").Value = "  -6.63%  "

# Row 39 - Kaspa
$ws.Range("E39").Value = "  +0.60%  "

# Row 40 - Stellar
$ws.Range("E40").Value = "  +0.53%  "

# Row 41 - RenderToken
$ws.Range("D41").Value = "'4.34"
$ws.Range("E41").Value = "  +9.12%  "

# Row 42 - EnergySwap
$ws.Range("D42").Value = "'22.15"
$ws.Range("E42").Value = "  -1.68%  "

# Row 43 - FirstDigitalUSD
$ws.Range("E43").Value = "  +0.06%  "

# Row 44 - NEARProtocol
$ws.Range("D44").Value = "'3.30"
$ws.Range("E44").Value = "  +1.67%  "

# Row 45 - VeChain
$ws.Range("E45").Value = "  -1.95%  "

# Row 46 - Maker
$ws.Range("D46").Value = "1.972.96"
$ws.Range("E46").Value = "  -1.33%  "

# Row 47 - FraxShare
$ws.Range("D47").Value = "'8.86"

# Row 48 - BitcoinSV
$ws.Range("D48").Value = "'81.07"
$ws.Range("E48").Value = "  -2.03%  "

# Row 49 - Algorand
$ws.Range("D49").Value = "'0.192"
$ws.Range("E49").Value = "  -0.86%  "

# Row 50 - SEI
$ws.Range("D50").Value = "'0.854"
$ws.Range("E50").Value = "  +10.41%  "

# Row 51 - RocketPoolETH
$ws.Range("D51").Value = "2.729.32"
$ws.Range("E51").Value = "  -2.92%  "
